# Auto-generated Excel COM-interop script
$wb = $excel.ActiveWorkbook

# Step 1: rename existing "总计" sheet (sheetId=5) to "2022-Q1" and make it the new fund table
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

# Step 2: add a brand-new sheet right after "2022-Q1" (this becomes the new "总计", sheetId=6)
$zj = $wb.Worksheets.Add($null, $q1)
$zj.Name = "总计"

# --- Build header + styles for "2022-Q1" (fund table, A1:H40) ---
$fmtHeader = $q1.Range("B1")   # pre-existing header cell already styled s=2
$fmtIndex  = $q1.Range("A2")   # pre-existing column-A index cell already styled s=2

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# copy the header style (bold/border/center-top) across B1:H1
$fmtHeader.Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

# data rows 2..40
$q1.Cells.Item(2, 1).Value = 0
$q1.Cells.Item(2, 2).Value = "'001417"
$q1.Cells.Item(2, 3).Value = "汇添富医疗服务灵活配置混合"
$q1.Cells.Item(2, 4).Value = "'38.45"
$q1.Cells.Item(2, 5).Value = "'77.97"
$q1.Cells.Item(2, 6).Value = "'6.51"
$q1.Cells.Item(2, 7).Value = "'2.5031"
$q1.Cells.Item(2, 8).Value = 2
$q1.Cells.Item(3, 1).Value = 1
$q1.Cells.Item(3, 2).Value = "'015122"
$q1.Cells.Item(3, 3).Value = "汇添富医疗服务灵活配置混合D"
$q1.Cells.Item(3, 4).Value = "'38.45"
$q1.Cells.Item(3, 5).Value = "'77.97"
$q1.Cells.Item(3, 6).Value = "'6.51"
$q1.Cells.Item(3, 7).Value = "'2.5031"
$q1.Cells.Item(3, 8).Value = 2
$q1.Cells.Item(4, 1).Value = 2
$q1.Cells.Item(4, 2).Value = "'519018"
$q1.Cells.Item(4, 3).Value = "汇添富均衡增长混合"
$q1.Cells.Item(4, 4).Value = "'39.45"
$q1.Cells.Item(4, 5).Value = "'93.81"
$q1.Cells.Item(4, 6).Value = "'5.09"
$q1.Cells.Item(4, 7).Value = "'2.0080"
$q1.Cells.Item(4, 8).Value = 4
$q1.Cells.Item(5, 1).Value = 3
$q1.Cells.Item(5, 2).Value = "'010481"
$q1.Cells.Item(5, 3).Value = "汇添富高质量成长精选2年持有期混合"
$q1.Cells.Item(5, 4).Value = "'36.73"
$q1.Cells.Item(5, 5).Value = "'72.35"
$q1.Cells.Item(5, 6).Value = "'4.17"
$q1.Cells.Item(5, 7).Value = "'1.5316"
$q1.Cells.Item(5, 8).Value = 5
$q1.Cells.Item(6, 1).Value = 4
$q1.Cells.Item(6, 2).Value = "'110023"
$q1.Cells.Item(6, 3).Value = "易方达医疗保健行业混合"
$q1.Cells.Item(6, 4).Value = "'36.81"
$q1.Cells.Item(6, 5).Value = "'90.49"
$q1.Cells.Item(6, 6).Value = "'4.12"
$q1.Cells.Item(6, 7).Value = "'1.5166"
$q1.Cells.Item(6, 8).Value = 6
$q1.Cells.Item(7, 1).Value = 5
$q1.Cells.Item(7, 2).Value = "'005453"
$q1.Cells.Item(7, 3).Value = "前海开源医疗健康灵活配置混合A"
$q1.Cells.Item(7, 4).Value = "'24.34"
$q1.Cells.Item(7, 5).Value = "'89.12"
$q1.Cells.Item(7, 6).Value = "'5.73"
$q1.Cells.Item(7, 7).Value = "'1.3947"
$q1.Cells.Item(7, 8).Value = 7
$q1.Cells.Item(8, 1).Value = 6
$q1.Cells.Item(8, 2).Value = "'010599"
$q1.Cells.Item(8, 3).Value = "汇添富高质量成长30一年持有期混合A"
$q1.Cells.Item(8, 4).Value = "'20.34"
$q1.Cells.Item(8, 5).Value = "'74.07"
$q1.Cells.Item(8, 6).Value = "'5.99"
$q1.Cells.Item(8, 7).Value = "'1.2184"
$q1.Cells.Item(8, 8).Value = 1
$q1.Cells.Item(9, 1).Value = 7
$q1.Cells.Item(9, 2).Value = "'012155"
$q1.Cells.Item(9, 3).Value = "汇添富成长先锋六个月持有期混合型证券投资基金A"
$q1.Cells.Item(9, 4).Value = "'21.01"
$q1.Cells.Item(9, 5).Value = "'71.54"
$q1.Cells.Item(9, 6).Value = "'5.61"
$q1.Cells.Item(9, 7).Value = "'1.1787"
$q1.Cells.Item(9, 8).Value = 2
$q1.Cells.Item(10, 1).Value = 8
$q1.Cells.Item(10, 2).Value = "'011271"
$q1.Cells.Item(10, 3).Value = "汇添富价值成长均衡投资混合A"
$q1.Cells.Item(10, 4).Value = "'24.78"
$q1.Cells.Item(10, 5).Value = "'92.58"
$q1.Cells.Item(10, 6).Value = "'3.50"
$q1.Cells.Item(10, 7).Value = "'0.8673"
$q1.Cells.Item(10, 8).Value = 10
$q1.Cells.Item(11, 1).Value = 9
$q1.Cells.Item(11, 2).Value = "'010387"
$q1.Cells.Item(11, 3).Value = "易方达医药生物股票A"
$q1.Cells.Item(11, 4).Value = "'22.51"
$q1.Cells.Item(11, 5).Value = "'83.49"
$q1.Cells.Item(11, 6).Value = "'3.39"
$q1.Cells.Item(11, 7).Value = "'0.7631"
$q1.Cells.Item(11, 8).Value = 9
$q1.Cells.Item(12, 1).Value = 10
$q1.Cells.Item(12, 2).Value = "'005454"
$q1.Cells.Item(12, 3).Value = "前海开源医疗健康灵活配置混合C"
$q1.Cells.Item(12, 4).Value = "'12.66"
$q1.Cells.Item(12, 5).Value = "'89.12"
$q1.Cells.Item(12, 6).Value = "'5.73"
$q1.Cells.Item(12, 7).Value = "'0.7254"
$q1.Cells.Item(12, 8).Value = 7
$q1.Cells.Item(13, 1).Value = 11
$q1.Cells.Item(13, 2).Value = "'000913"
$q1.Cells.Item(13, 3).Value = "农银医疗保健主题股票"
$q1.Cells.Item(13, 4).Value = "'22.47"
$q1.Cells.Item(13, 5).Value = "'92.80"
$q1.Cells.Item(13, 6).Value = "'2.95"
$q1.Cells.Item(13, 7).Value = "'0.6629"
$q1.Cells.Item(13, 8).Value = 9
$q1.Cells.Item(14, 1).Value = 12
$q1.Cells.Item(14, 2).Value = "'470008"
$q1.Cells.Item(14, 3).Value = "汇添富策略回报混合"
$q1.Cells.Item(14, 4).Value = "'12.75"
$q1.Cells.Item(14, 5).Value = "'94.01"
$q1.Cells.Item(14, 6).Value = "'4.91"
$q1.Cells.Item(14, 7).Value = "'0.6260"
$q1.Cells.Item(14, 8).Value = 4
$q1.Cells.Item(15, 1).Value = 13
$q1.Cells.Item(15, 2).Value = "'008293"
$q1.Cells.Item(15, 3).Value = "农银汇理创新医疗混合"
$q1.Cells.Item(15, 4).Value = "'10.12"
$q1.Cells.Item(15, 5).Value = "'91.66"
$q1.Cells.Item(15, 6).Value = "'2.96"
$q1.Cells.Item(15, 7).Value = "'0.2996"
$q1.Cells.Item(15, 8).Value = 9
$q1.Cells.Item(16, 1).Value = 14
$q1.Cells.Item(16, 2).Value = "'001766"
$q1.Cells.Item(16, 3).Value = "上投摩根医疗健康股票"
$q1.Cells.Item(16, 4).Value = "'10.35"
$q1.Cells.Item(16, 5).Value = "'80.54"
$q1.Cells.Item(16, 6).Value = "'2.71"
$q1.Cells.Item(16, 7).Value = "'0.2805"
$q1.Cells.Item(16, 8).Value = 9
$q1.Cells.Item(17, 1).Value = 15
$q1.Cells.Item(17, 2).Value = "'009468"
$q1.Cells.Item(17, 3).Value = "博时健康成长主题双周定期可赎回混合A"
$q1.Cells.Item(17, 4).Value = "'6.39"
$q1.Cells.Item(17, 5).Value = "'85.57"
$q1.Cells.Item(17, 6).Value = "'3.52"
$q1.Cells.Item(17, 7).Value = "'0.2249"
$q1.Cells.Item(17, 8).Value = 6
$q1.Cells.Item(18, 1).Value = 16
$q1.Cells.Item(18, 2).Value = "'090020"
$q1.Cells.Item(18, 3).Value = "大成健康产业混合"
$q1.Cells.Item(18, 4).Value = "'3.76"
$q1.Cells.Item(18, 5).Value = "'91.73"
$q1.Cells.Item(18, 6).Value = "'5.18"
$q1.Cells.Item(18, 7).Value = "'0.1948"
$q1.Cells.Item(18, 8).Value = 10
$q1.Cells.Item(19, 1).Value = 17
$q1.Cells.Item(19, 2).Value = "'001898"
$q1.Cells.Item(19, 3).Value = "易方达大健康主题灵活配置混合"
$q1.Cells.Item(19, 4).Value = "'4.33"
$q1.Cells.Item(19, 5).Value = "'87.42"
$q1.Cells.Item(19, 6).Value = "'4.29"
$q1.Cells.Item(19, 7).Value = "'0.1858"
$q1.Cells.Item(19, 8).Value = 6
$q1.Cells.Item(20, 1).Value = 18
$q1.Cells.Item(20, 2).Value = "'010388"
$q1.Cells.Item(20, 3).Value = "易方达医药生物股票C"
$q1.Cells.Item(20, 4).Value = "'4.78"
$q1.Cells.Item(20, 5).Value = "'83.49"
$q1.Cells.Item(20, 6).Value = "'3.39"
$q1.Cells.Item(20, 7).Value = "'0.1620"
$q1.Cells.Item(20, 8).Value = 9
$q1.Cells.Item(21, 1).Value = 19
$q1.Cells.Item(21, 2).Value = "'012045"
$q1.Cells.Item(21, 3).Value = "大成医药健康股票A"
$q1.Cells.Item(21, 4).Value = "'2.87"
$q1.Cells.Item(21, 5).Value = "'93.58"
$q1.Cells.Item(21, 6).Value = "'4.83"
$q1.Cells.Item(21, 7).Value = "'0.1386"
$q1.Cells.Item(21, 8).Value = 10
$q1.Cells.Item(22, 1).Value = 20
$q1.Cells.Item(22, 2).Value = "'011601"
$q1.Cells.Item(22, 3).Value = "前海开源公共卫生主题精选股票A"
$q1.Cells.Item(22, 4).Value = "'2.25"
$q1.Cells.Item(22, 5).Value = "'88.30"
$q1.Cells.Item(22, 6).Value = "'5.76"
$q1.Cells.Item(22, 7).Value = "'0.1296"
$q1.Cells.Item(22, 8).Value = 7
$q1.Cells.Item(23, 1).Value = 21
$q1.Cells.Item(23, 2).Value = "'011288"
$q1.Cells.Item(23, 3).Value = "上银医疗健康混合A"
$q1.Cells.Item(23, 4).Value = "'1.48"
$q1.Cells.Item(23, 5).Value = "'81.70"
$q1.Cells.Item(23, 6).Value = "'4.78"
$q1.Cells.Item(23, 7).Value = "'0.0707"
$q1.Cells.Item(23, 8).Value = 4
$q1.Cells.Item(24, 1).Value = 22
$q1.Cells.Item(24, 2).Value = "'011259"
$q1.Cells.Item(24, 3).Value = "汇添富高质量成长30一年持有期混合C"
$q1.Cells.Item(24, 4).Value = "'0.90"
$q1.Cells.Item(24, 5).Value = "'74.07"
$q1.Cells.Item(24, 6).Value = "'5.99"
$q1.Cells.Item(24, 7).Value = "'0.0539"
$q1.Cells.Item(24, 8).Value = 1
$q1.Cells.Item(25, 1).Value = 23
$q1.Cells.Item(25, 2).Value = "'011602"
$q1.Cells.Item(25, 3).Value = "前海开源公共卫生主题精选股票C"
$q1.Cells.Item(25, 4).Value = "'0.86"
$q1.Cells.Item(25, 5).Value = "'88.30"
$q1.Cells.Item(25, 6).Value = "'5.76"
$q1.Cells.Item(25, 7).Value = "'0.0495"
$q1.Cells.Item(25, 8).Value = 7
$q1.Cells.Item(26, 1).Value = 24
$q1.Cells.Item(26, 2).Value = "'009469"
$q1.Cells.Item(26, 3).Value = "博时健康成长主题双周定期可赎回混合C"
$q1.Cells.Item(26, 4).Value = "'1.11"
$q1.Cells.Item(26, 5).Value = "'85.57"
$q1.Cells.Item(26, 6).Value = "'3.52"
$q1.Cells.Item(26, 7).Value = "'0.0391"
$q1.Cells.Item(26, 8).Value = 6
$q1.Cells.Item(27, 1).Value = 25
$q1.Cells.Item(27, 2).Value = "'011765"
$q1.Cells.Item(27, 3).Value = "兴银高端制造混合A"
$q1.Cells.Item(27, 4).Value = "'1.01"
$q1.Cells.Item(27, 5).Value = "'93.23"
$q1.Cells.Item(27, 6).Value = "'3.50"
$q1.Cells.Item(27, 7).Value = "'0.0354"
$q1.Cells.Item(27, 8).Value = 1
$q1.Cells.Item(28, 1).Value = 26
$q1.Cells.Item(28, 2).Value = "'002863"
$q1.Cells.Item(28, 3).Value = "金信深圳成长灵活配置混合"
$q1.Cells.Item(28, 4).Value = "'0.44"
$q1.Cells.Item(28, 5).Value = "'94.54"
$q1.Cells.Item(28, 6).Value = "'7.87"
$q1.Cells.Item(28, 7).Value = "'0.0346"
$q1.Cells.Item(28, 8).Value = 2
$q1.Cells.Item(29, 1).Value = 27
$q1.Cells.Item(29, 2).Value = "'013441"
$q1.Cells.Item(29, 3).Value = "西藏东财创新医疗六个月定开混合"
$q1.Cells.Item(29, 4).Value = "'0.58"
$q1.Cells.Item(29, 5).Value = "'81.46"
$q1.Cells.Item(29, 6).Value = "'5.70"
$q1.Cells.Item(29, 7).Value = "'0.0331"
$q1.Cells.Item(29, 8).Value = 6
$q1.Cells.Item(30, 1).Value = 28
$q1.Cells.Item(30, 2).Value = "'012156"
$q1.Cells.Item(30, 3).Value = "汇添富成长先锋六个月持有期混合型证券投资基金C"
$q1.Cells.Item(30, 4).Value = "'0.57"
$q1.Cells.Item(30, 5).Value = "'71.54"
$q1.Cells.Item(30, 6).Value = "'5.61"
$q1.Cells.Item(30, 7).Value = "'0.0320"
$q1.Cells.Item(30, 8).Value = 2
$q1.Cells.Item(31, 1).Value = 29
$q1.Cells.Item(31, 2).Value = "'217021"
$q1.Cells.Item(31, 3).Value = "招商优势企业混合"
$q1.Cells.Item(31, 4).Value = "'0.36"
$q1.Cells.Item(31, 5).Value = "'69.72"
$q1.Cells.Item(31, 6).Value = "'5.15"
$q1.Cells.Item(31, 7).Value = "'0.0185"
$q1.Cells.Item(31, 8).Value = 5
$q1.Cells.Item(32, 1).Value = 30
$q1.Cells.Item(32, 2).Value = "'011766"
$q1.Cells.Item(32, 3).Value = "兴银高端制造混合C"
$q1.Cells.Item(32, 4).Value = "'0.39"
$q1.Cells.Item(32, 5).Value = "'93.23"
$q1.Cells.Item(32, 6).Value = "'3.50"
$q1.Cells.Item(32, 7).Value = "'0.0136"
$q1.Cells.Item(32, 8).Value = 1
$q1.Cells.Item(33, 1).Value = 31
$q1.Cells.Item(33, 2).Value = "'012046"
$q1.Cells.Item(33, 3).Value = "大成医药健康股票C"
$q1.Cells.Item(33, 4).Value = "'0.25"
$q1.Cells.Item(33, 5).Value = "'93.58"
$q1.Cells.Item(33, 6).Value = "'4.83"
$q1.Cells.Item(33, 7).Value = "'0.0121"
$q1.Cells.Item(33, 8).Value = 10
$q1.Cells.Item(34, 1).Value = 32
$q1.Cells.Item(34, 2).Value = "'011272"
$q1.Cells.Item(34, 3).Value = "汇添富价值成长均衡投资混合C"
$q1.Cells.Item(34, 4).Value = "'0.32"
$q1.Cells.Item(34, 5).Value = "'92.58"
$q1.Cells.Item(34, 6).Value = "'3.50"
$q1.Cells.Item(34, 7).Value = "'0.0112"
$q1.Cells.Item(34, 8).Value = 10
$q1.Cells.Item(35, 1).Value = 33
$q1.Cells.Item(35, 2).Value = "'011289"
$q1.Cells.Item(35, 3).Value = "上银医疗健康混合C"
$q1.Cells.Item(35, 4).Value = "'0.23"
$q1.Cells.Item(35, 5).Value = "'81.70"
$q1.Cells.Item(35, 6).Value = "'4.78"
$q1.Cells.Item(35, 7).Value = "'0.0110"
$q1.Cells.Item(35, 8).Value = 4
$q1.Cells.Item(36, 1).Value = 34
$q1.Cells.Item(36, 2).Value = "'008037"
$q1.Cells.Item(36, 3).Value = "兴银先锋成长混合A"
$q1.Cells.Item(36, 4).Value = "'0.41"
$q1.Cells.Item(36, 5).Value = "'79.32"
$q1.Cells.Item(36, 6).Value = "'2.38"
$q1.Cells.Item(36, 7).Value = "'0.0098"
$q1.Cells.Item(36, 8).Value = 3
$q1.Cells.Item(37, 1).Value = 35
$q1.Cells.Item(37, 2).Value = "'006240"
$q1.Cells.Item(37, 3).Value = "中融医疗健康精选混合A"
$q1.Cells.Item(37, 4).Value = "'0.14"
$q1.Cells.Item(37, 5).Value = "'94.05"
$q1.Cells.Item(37, 6).Value = "'3.31"
$q1.Cells.Item(37, 7).Value = "'0.0046"
$q1.Cells.Item(37, 8).Value = 10
$q1.Cells.Item(38, 1).Value = 36
$q1.Cells.Item(38, 2).Value = "'008038"
$q1.Cells.Item(38, 3).Value = "兴银先锋成长混合C"
$q1.Cells.Item(38, 4).Value = "'0.17"
$q1.Cells.Item(38, 5).Value = "'79.32"
$q1.Cells.Item(38, 6).Value = "'2.38"
$q1.Cells.Item(38, 7).Value = "'0.0040"
$q1.Cells.Item(38, 8).Value = 3
$q1.Cells.Item(39, 1).Value = 37
$q1.Cells.Item(39, 2).Value = "'006241"
$q1.Cells.Item(39, 3).Value = "中融医疗健康精选混合C"
$q1.Cells.Item(39, 4).Value = "'0.08"
$q1.Cells.Item(39, 5).Value = "'94.05"
$q1.Cells.Item(39, 6).Value = "'3.31"
$q1.Cells.Item(39, 7).Value = "'0.0026"
$q1.Cells.Item(39, 8).Value = 10
$q1.Cells.Item(40, 1).Value = 38
$q1.Cells.Item(40, 2).Value = "'005146"
$q1.Cells.Item(40, 3).Value = "兴银丰润灵活配置混合"
$q1.Cells.Item(40, 4).Value = "'0.05"
$q1.Cells.Item(40, 5).Value = "'93.36"
$q1.Cells.Item(40, 6).Value = "'4.10"
$q1.Cells.Item(40, 7).Value = "'0.0020"
$q1.Cells.Item(40, 8).Value = 1

# copy column-A index style (bold/border/center-top) down A2:A40
$fmtIndex.Copy()
$q1.Range("A2:A40").PasteSpecial(-4122)

# --- Build header + styles for the new "总计" sheet (A1:D6) ---
$zj.Range("B1").Value = "日期"
$zj.Range("C1").Value = "持有数量(只)"
$zj.Range("D1").Value = "持有市值(亿元)"

$fmtHeader.Copy()
$zj.Range("B1:D1").PasteSpecial(-4122)

# data rows 2..6
$zj.Cells.Item(2, 1).Value = 0
$zj.Cells.Item(2, 2).Value = "2022-Q1"
$zj.Cells.Item(2, 3).Value = 39
$zj.Cells.Item(2, 4).Value = 19.55
$zj.Cells.Item(3, 1).Value = 1
$zj.Cells.Item(3, 2).Value = "2021-Q4"
$zj.Cells.Item(3, 3).Value = 68
$zj.Cells.Item(3, 4).Value = 31.08
$zj.Cells.Item(4, 1).Value = 2
$zj.Cells.Item(4, 2).Value = "2021-Q3"
$zj.Cells.Item(4, 3).Value = 66
$zj.Cells.Item(4, 4).Value = 29.9
$zj.Cells.Item(5, 1).Value = 3
$zj.Cells.Item(5, 2).Value = "2021-Q2"
$zj.Cells.Item(5, 3).Value = 38
$zj.Cells.Item(5, 4).Value = 12.69
$zj.Cells.Item(6, 1).Value = 4
$zj.Cells.Item(6, 2).Value = "2021-Q1"
$zj.Cells.Item(6, 3).Value = 6
$zj.Cells.Item(6, 4).Value = 0.49

$fmtIndex.Copy()
$zj.Range("A2:A6").PasteSpecial(-4122)

Write-Output "done"
